$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"
$ws.Range("B2").Select()
